$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Glayds  Bundotich", shifting the rows below it up by one
$ws.Rows("3:3").Delete()

# Keep the numeric-looking figures stored as plain text (as in the source data)
$ws.Range("B2:E6").NumberFormat = "@"

# Row 2 - Beatrice Chege
$ws.Range("A2").Value = "Beatrice Chege"
$ws.Range("B2").Value = "7.00"
$ws.Range("C2").Value = "50.00"
$ws.Range("D2").Value = "-43.00"
$ws.Range("E2").Value = "14.00%"

# Row 3 - Jane Gichohi
$ws.Range("A3").Value = "Jane Gichohi"
$ws.Range("B3").Value = "5.00"
$ws.Range("C3").Value = "30.00"
$ws.Range("D3").Value = "-25.00"
$ws.Range("E3").Value = "16.67%"

# Row 4 - Mirriam Makau
$ws.Range("A4").Value = "Mirriam Makau"
$ws.Range("B4").Value = "1.00"
$ws.Range("C4").Value = "25.00"
$ws.Range("D4").Value = "-24.00"
$ws.Range("E4").Value = "4.00%"

# Row 5 - Victor Njogu
$ws.Range("A5").Value = "Victor Njogu"
$ws.Range("B5").Value = "1.00"
$ws.Range("C5").Value = "40.00"
$ws.Range("D5").Value = "-39.00"
$ws.Range("E5").Value = "2.50%"

# Row 6 - KD Totals
$ws.Range("A6").Value = "KD Totals"
$ws.Range("B6").Value = "14.00"
$ws.Range("C6").Value = "145.00"
$ws.Range("D6").Value = "-131.00"
$ws.Range("E6").Value = "37.17%"
